# heroes.xlsx update - "Added the fields placeholder"
#
# Rebalances the B-column (exp-per-level) growth curve on the four
# "*_levels" sheets (ulko_levels, ferlin_levels, phoebey_levels,
# bob_levels) from a doubling curve (B*2) to a gentler
# ROUND(prev*1.2+k, 0) curve, and moves the active selection/active
# sheet around (cosmetic UI state) to match the author's last-saved
# view.

$wb = $excel.ActiveWorkbook

$ulko    = $wb.Worksheets.Item("ulko_levels")
$ferlin  = $wb.Worksheets.Item("ferlin_levels")
$phoebey = $wb.Worksheets.Item("phoebey_levels")
$bob     = $wb.Worksheets.Item("bob_levels")
$heroes  = $wb.Worksheets.Item("heroes")

# --- Rebalance the B column (same new curve on all four sheets) ---
foreach ($ws in @($ulko, $ferlin, $phoebey, $bob)) {
    $ws.Range("B3").Value = 6
    $ws.Range("B4").Formula  = "=ROUND(B3*1.2+6, 0)"
    $ws.Range("B5").Formula  = "=ROUND(B4*1.2+4, 0)"
    $ws.Range("B6").Formula  = "=ROUND(B5*1.2+4, 0)"
    $ws.Range("B7").Formula  = "=ROUND(B6*1.2+4, 0)"
    $ws.Range("B8").Formula  = "=ROUND(B7*1.2+4, 0)"
    $ws.Range("B9").Formula  = "=ROUND(B8*1.2+4, 0)"
    $ws.Range("B10").Formula = "=ROUND(B9*1.2+4, 0)"
    $ws.Range("B11").Formula = "=ROUND(B10*1.2+3, 0)"
}

# --- Per-sheet selection (cursor position) update ---
[void]$heroes.Range("C10").Select()
[void]$ulko.Range("B18").Select()
[void]$ferlin.Range("C17").Select()
[void]$phoebey.Range("D23").Select()
[void]$bob.Range("C19").Select()

# --- Move the active tab from "heroes" to "bob_levels" ---
[void]$bob.Activate()

# --- Best-effort: restore the saved window position/size (cosmetic;
#     not guaranteed to round-trip through every host, but harmless) ---
$excel.ActiveWindow.Left = 1395
$excel.ActiveWindow.Top = 3465
$excel.ActiveWindow.Width = 33915
$excel.ActiveWindow.Height = 15150
